# Applies the edit described by the diff:
#  - Sheet "s__Mesosutterella massiliensis-b-p": delete row 4
#    (label_GCF_001182045_3.fasta), shifting rows 5-8 up to 4-7.
#  - Sheet "s__Mesosutterella multiformis-b-p": delete rows 8-15
#    (UMGS117_2, _22, _3, _30, _36, _38, _4, _6), shifting rows 16-46
#    up to 8-38.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("s__Mesosutterella massiliensis-b-p")
$ws1.Rows.Item(4).Delete()

$ws2 = $wb.Worksheets.Item("s__Mesosutterella multiformis-b-p")
$ws2.Range("A8:A15").EntireRow.Delete()
